$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values (NN correction) for rows 2-9, columns B,C,D,E,F,H
# Column A and G are unchanged.

$data = @{
    2 = @{ B = 1314.430908203125;  C = 0.9861;             D = 0.9617000222206116; E = 1.40719997882843;  F = 0.8133000135421753; H = 1.1255 }
    3 = @{ B = 1182.26953125;      C = 0.9458;              D = 0.9073;             E = 1.34060001373291;  F = 0.800599992275238;  H = 0.6443 }
    4 = @{ B = 810.4566040039062;  C = 0.9468;              D = 0.9135;             E = 1.288699984550476; F = 0.8578000068664551; H = 0.6992 }
    5 = @{ B = 880.3624267578125;  C = 0.9287;              D = 0.9183;             E = 1.080000042915344; F = 0.7282999753952026; H = 0.7412 }
    6 = @{ B = 1165.527709960938;  C = 0.9221;              D = 0.9177;             E = 1.07260000705719;  F = 0.8080999851226807; H = 0.7363 }
    7 = @{ B = 908.8593139648438;  C = 0.9153;              D = 0.9110000133514404; E = 0.9929999709129333;F = 0.8610000014305115; H = 0.6768999999999999 }
    8 = @{ B = 1010.140930175781;  C = 0.9051;              D = 0.9036999999999999;E = 1.039399981498718; F = 0.8550999760627747; H = 0.6118 }
    9 = @{ B = 7272.0478515625;    C = 0.9372;              D = 0.9121;             E = 1.40719997882843;  F = 0.7282999753952026; H = 5.2352 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("H$row").Value = $vals.H
}
